$d = $word.ActiveDocument

$pairs = @(
    @("45×43=1935", "26×67=1742"),
    @("25×29=725", "26×69=1794"),
    @("86×71=6106", "27×53=1431"),
    @("43×71=3053", "44×17=748"),
    @("87×84=7308", "66×73=4818"),
    @("55×90=4950", "54×85=4590"),
    @("46×58=2668", "31×14=434"),
    @("93×39=3627", "30×17=510"),
    @("84×84=7056", "82×87=7134"),
    @("62×22=1364", "24×48=1152"),
    @("90×80=7200", "60×19=1140"),
    @("36×25=900", "32×82=2624"),
    @("63×26=1638", "40×80=3200"),
    @("95×37=3515", "36×28=1008"),
    @("58×98=5684", "53×27=1431"),
    @("52×14=728", "54×88=4752"),
    @("58×19=1102", "18×44=792"),
    @("68×62=4216", "27×94=2538"),
    @("48×30=1440", "42×77=3234"),
    @("79×94=7426", "56×37=2072"),
    @("34×42=1428", "35×48=1680"),
    @("61×29=1769", "25×57=1425"),
    @("93×49=4557", "25×43=1075"),
    @("39×17=663", "74×96=7104"),
    @("96×54=5184", "82×65=5330")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
